# Feria Lagunitas de Puerto Montt - Apio
# A new weekly price-report row is inserted right after the existing
# row 127 (pushing the existing rows 128-196 down to 129-197), which is
# why the sheet's used range grows from A1:R196 to A1:R197.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at position 128; everything from the old
# row 128 onward shifts down by one row (old 196 becomes new 197).
$ws.Rows.Item(128).Insert()

# Populate the newly inserted row 128 with the new weekly record.
$ws.Range("A128").Value = 4
$ws.Range("B128").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C128").Value = "Los Lagos"
$ws.Range("D128").Value = 44596
$ws.Range("E128").Value = 10
$ws.Range("F128").Value = 100112017
$ws.Range("G128").Value = "Apio"
$ws.Range("H128").Value = "Americana (o)"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 40
$ws.Range("K128").Value = 11000
$ws.Range("L128").Value = 11000
$ws.Range("M128").Value = 11000
$ws.Range("N128").Value = "$/docena de matas"
$ws.Range("O128").Value = "Región de Coquimbo"
$ws.Range("P128").Value = 1833
$ws.Range("Q128").Value = 6
$ws.Range("R128").Value = "Hortaliza"
